$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string run text) ---
$ws.Range("A8").Value = "Volume 30   Number  1"
$ws.Range("C9").Value = "Report Covering the Week  1/2/2023  Through  1/8/2023"

# --- Data table cell updates (rows 16-29, 38-43) ---
# Row 16
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -50
$ws.Range("I16").NumberFormat = '#,##0'
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = -25
$ws.Range("L16").Value = 50
$ws.Range("M16").Value = -57.142857142857
$ws.Range("N16").Value = -87.5
# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 16.666666666666
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 21.428571428571
$ws.Range("I17").NumberFormat = '#,##0'
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 8
$ws.Range("K17").Value = -12.5
$ws.Range("L17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L17").Value = 600
$ws.Range("M17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M17").Value = 600
$ws.Range("N17").Value = -46.153846153846
# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -69.230769230769
$ws.Range("I18").NumberFormat = '#,##0'
$ws.Range("I18").Value = 2
$ws.Range("J18").NumberFormat = '#,##0'
$ws.Range("J18").Value = 7
$ws.Range("K18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K18").Value = -71.428571428571
$ws.Range("L18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L18").Value = -50
$ws.Range("N18").Value = -88.235294117647
# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -14.285714285714
$ws.Range("I19").NumberFormat = '#,##0'
$ws.Range("I19").Value = 6
$ws.Range("J19").Value = 13
$ws.Range("K19").Value = -53.846153846153
$ws.Range("L19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L19").Value = -40
$ws.Range("M19").Value = -45.454545454545
$ws.Range("N19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N19").Value = -75
# Row 20
$ws.Range("C20").Value = 1
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "0"
$ws.Range("G20").NumberFormat = "General"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "***.*"
$ws.Range("H20").NumberFormat = "General"
$ws.Range("I20").NumberFormat = '#,##0'
$ws.Range("I20").Value = 1
$ws.Range("M20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -92.307692307692
# Row 21
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -32.142857142857
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -23.300970873786
$ws.Range("I21").NumberFormat = '#,##0'
$ws.Range("I21").Value = 19
$ws.Range("J21").Value = 32
$ws.Range("K21").Value = -40.625
$ws.Range("L21").Value = 11.764705882352
$ws.Range("M21").Value = -5
$ws.Range("N21").Value = -79.120879120879
# Row 22
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").NumberFormat = '#,##0'
$ws.Range("I22").Value = 1
$ws.Range("M22").Value = -75
# Row 23
$ws.Range("C23").NumberFormat = '#,##0'
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -77.777777777777
$ws.Range("I23").NumberFormat = '#,##0'
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 3
$ws.Range("K23").Value = -66.666666666666
# Row 24
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 63.157894736842
$ws.Range("F24").Value = 128
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = 48.837209302325
$ws.Range("I24").Value = 32
$ws.Range("J24").NumberFormat = '#,##0'
$ws.Range("J24").Value = 19
$ws.Range("K24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K24").Value = 68.421052631578
$ws.Range("L24").Value = 68.421052631578
$ws.Range("M24").Value = 10.344827586206
# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 266.666666666667
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 113.333333333333
$ws.Range("I25").Value = 13
$ws.Range("J25").NumberFormat = '#,##0'
$ws.Range("J25").Value = 3
$ws.Range("K25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K25").Value = 333.333333333333
$ws.Range("L25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L25").Value = 333.333333333333
$ws.Range("M25").Value = 44.444444444444
# Row 26
$ws.Range("L26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L26").Value = -100
# Row 27
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("I27").NumberFormat = '#,##0'
$ws.Range("I27").Value = 1
$ws.Range("J27").NumberFormat = '#,##0'
$ws.Range("J27").Value = 1
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K27").Value = 0
$ws.Range("L27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L27").Value = -50
# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("G28").NumberFormat = "General"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "***.*"
$ws.Range("H28").NumberFormat = "General"
# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"
$ws.Range("G29").NumberFormat = "General"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "***.*"
$ws.Range("H29").NumberFormat = "General"
# Row 38
$ws.Range("J38").Value = 136
$ws.Range("K38").Value = -51.254480286738
$ws.Range("L38").Value = -70.562770562770
$ws.Range("M38").Value = -89.649923896499
$ws.Range("N38").Value = -92.872117400419
# Row 40
$ws.Range("J40").Value = 178
$ws.Range("K40").Value = -42.948717948717
$ws.Range("L40").Value = -54.591836734693
$ws.Range("M40").Value = -78.292682926829
$ws.Range("N40").Value = -88.611644273832
# Row 41
$ws.Range("J41").Value = 636
$ws.Range("K41").Value = -13.586956521739
$ws.Range("L41").Value = -13.114754098360
$ws.Range("M41").Value = -44.014084507042
$ws.Range("N41").Value = -62.915451895043
# Row 43
$ws.Range("J43").Value = 1207
$ws.Range("K43").Value = -28.325415676959
$ws.Range("L43").Value = -44.016697588126
$ws.Range("M43").Value = -72.985675917636
$ws.Range("N43").Value = -81.530221882172
